$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: AH7 1 -> 0, AP7 13 -> 12 ---
$ws.Range("AH7").Value = 0
$ws.Range("AP7").Value = 12

# --- AP column updates for rows 24-221 (recomputed weighted scores) ---
$ws.Range("AP24").Value = 0.30555555555
$ws.Range("AP25").Value = 0.36805555555
$ws.Range("AP26").Value = 0.4513888888833333
$ws.Range("AP27:AP28").Value = 0.5505952381083331
$ws.Range("AP29:AP51").Value = 0.8333333333416665
$ws.Range("AP52:AP71").Value = 0.9166666666749999
$ws.Range("AP72:AP91").Value = 0.9062500000083332
$ws.Range("AP92:AP109").Value = 0.8750000000083332
$ws.Range("AP110:AP115").Value = 0.7916666666749997
$ws.Range("AP116:AP122").Value = 0.7291666666749999
$ws.Range("AP123:AP129").Value = 0.63541666665
$ws.Range("AP130:AP132").Value = 0.5937499999833332
$ws.Range("AP133:AP221").Value = 0.340277777775

# --- New rows 222-233 (policy dates 9/30/2020 - 10/11/2020) ---
# Each new row duplicates the data pattern of row 221, with a new date label in column A
# and the recomputed AP score.
$newDates = @("9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020")
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 222 + $i
    $srcRow = $ws.Range("A221:AP221")
    $dstRow = $ws.Range("A" + $r + ":AP" + $r)
    $srcRow.Copy($dstRow)

    $a = $ws.Range("A" + $r)
    $a.NumberFormat = "@"
    $a.Value = $newDates[$i]
    $ws.Range("A221").Copy()
    $a.PasteSpecial(-4122)

    $ws.Range("AP" + $r).Value = 0.340277777775
}

Write-Host "edit complete"